# Fix text margins: zero out the text-box insets (lIns/tIns/bIns/rIns) on
# every shape's text frame across all slides, matching the upstream fix
# "text margins in pptx fixed" (a:bodyPr rtlCol="0" anchor="ctr" -> adds
# lIns="0" tIns="0" bIns="0" rIns="0").

$p = $ppt.ActivePresentation

for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $s = $p.Slides.Item($si)
    for ($i = 1; $i -le $s.Shapes.Count; $i++) {
        $sh = $s.Shapes.Item($i)
        if ($sh.HasTextFrame) {
            $tf = $sh.TextFrame
            $tf.MarginLeft = 0
            $tf.MarginTop = 0
            $tf.MarginBottom = 0
            $tf.MarginRight = 0
        }
    }
}
